{"js": "// Replace the 25 division-problem values in the single 5-column table.\n// The table has 20 rows total, but only every 4th row (0, 4, 8, 12, 16)\n// actually holds the \"NN\u00f7N=\" text; the rows in between are blank filler\n// rows used for students to write their answers. We update cell-by-cell\n// using row/column coordinates (not global text search) so that values\n// which coincidentally collide with other cells' old/new text never get\n// double-replaced.\n\nconst dataRowIndexes = [0, 4, 8, 12, 16];\n\n// Old -> new text, grouped by the same 5 \"data\" rows / 5 columns as the\n// source document, in reading order (row major).\nconst newValues = [\n  [\"55\u00f74=\", \"67\u00f75=\", \"28\u00f77=\", \"85\u00f78=\", \"96\u00f73=\"],\n  [\"27\u00f72=\", \"12\u00f77=\", \"82\u00f72=\", \"93\u00f77=\", \"36\u00f78=\"],\n  [\"63\u00f77=\", \"41\u00f73=\", \"70\u00f72=\", \"60\u00f77=\", \"14\u00f77=\"],\n  [\"61\u00f72=\", \"60\u00f77=\", \"28\u00f76=\", \"95\u00f75=\", \"45\u00f74=\"],\n  [\"57\u00f78=\", \"79\u00f76=\", \"78\u00f73=\", \"50\u00f77=\", \"87\u00f73=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (let r = 0; r < dataRowIndexes.length; r++) {\n  const tableRow = dataRowIndexes[r];\n  for (let c = 0; c < newValues[r].length; c++) {\n    const cell = table.getCell(tableRow, c);\n    // Replace just the text of the cell's range so the existing run\n    // formatting (font, size) and paragraph formatting (alignment) are\n    // preserved instead of being wiped out by a body-level insertText.\n    const range = cell.body.getRange();\n    range.insertText(newValues[r][c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 division-problem values in the single 5-column table.\n# The table has 20 rows total (1-based in the COM object model), but only\n# every 4th row (1, 5, 9, 13, 17) actually holds the \"NN\u00f7N=\" text; the\n# rows in between are blank filler rows used for students to write their\n# answers. We update cell-by-cell using row/column coordinates (not\n# global Find/Replace) so that values which coincidentally collide with\n# other cells' old/new text are never double-replaced.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$dataRows = @(1, 5, 9, 13, 17)\n\n$newValues = @(\n    @(\"55\u00f74=\", \"67\u00f75=\", \"28\u00f77=\", \"85\u00f78=\", \"96\u00f73=\"),\n    @(\"27\u00f72=\", \"12\u00f77=\", \"82\u00f72=\", \"93\u00f77=\", \"36\u00f78=\"),\n    @(\"63\u00f77=\", \"41\u00f73=\", \"70\u00f72=\", \"60\u00f77=\", \"14\u00f77=\"),\n    @(\"61\u00f72=\", \"60\u00f77=\", \"28\u00f76=\", \"95\u00f75=\", \"45\u00f74=\"),\n    @(\"57\u00f78=\", \"79\u00f76=\", \"78\u00f73=\", \"50\u00f77=\", \"87\u00f73=\")\n)\n\nfor ($r = 0; $r -lt $dataRows.Length; $r++) {\n    $rowIndex = $dataRows[$r]\n    $rowValues = $newValues[$r]\n    for ($c = 0; $c -lt $rowValues.Length; $c++) {\n        $cell = $tbl.Cell($rowIndex, $c + 1)\n        # Assigning Range.Text replaces just the cell's text run while\n        # keeping the existing run/paragraph formatting (font, size,\n        # alignment) intact.\n        $cell.Range.Text = $rowValues[$c]\n    }\n}\n"}
